# Update the Store_Data workbook:
#  - Credentials sheet: reset selection to a single cell (B2)
#  - Store_Details sheet: add new store-config columns (markup reorder +
#    pay_on_limit / order_approval / show_price_to_customer / store_fields /
#    manage_invoice) and refill all header + data cells, then reposition
#    the active selection.

$wb = $excel.ActiveWorkbook

# --- Credentials sheet -----------------------------------------------------
$wsCred = $wb.Worksheets.Item("Credentials")
[void]$wsCred.Select()
[void]$wsCred.Range("B2").Select()

# --- Store_Details sheet ----------------------------------------------------
$ws = $wb.Worksheets.Item("Store_Details")
[void]$ws.Select()

$ws.Range("A1").Value = "store_name"
$ws.Range("B1").Value = "email"
$ws.Range("C1").Value = "url_type"
$ws.Range("D1").Value = "subdomain_name"
$ws.Range("E1").Value = "main_store_url"
$ws.Range("F1").Value = "username"
$ws.Range("G1").Value = "password"
$ws.Range("H1").Value = "phone_number"
$ws.Range("I1").Value = "sales_agent"
$ws.Range("J1").Value = "add_as_front_customer"
$ws.Range("K1").Value = "markup_type"
$ws.Range("L1").Value = "markup_master"
$ws.Range("M1").Value = "flat_markup"
$ws.Range("N1").Value = "pay_on_account"
$ws.Range("O1").Value = "pay_on_limit"
$ws.Range("P1").Value = "open_b2b_store"
$ws.Range("Q1").Value = "department"
$ws.Range("R1").Value = "fix_billing_address"
$ws.Range("S1").Value = "fix_shipping_address"
$ws.Range("T1").Value = "order_approval"
$ws.Range("U1").Value = "show_price_to_customer"
$ws.Range("V1").Value = "quick_checkout"
$ws.Range("W1").Value = "store_fields"
$ws.Range("X1").Value = "allow_tax_exemption"
$ws.Range("Y1").Value = "manage_invoice"
$ws.Range("Z1").Value = "allow_partial_payment"
$ws.Range("AA1").Value = "notify"
$ws.Range("AB1").Value = "status"
$ws.Range("A2").Value = "09Prints"
$ws.Range("B2").Value = "vaibhav.darji+store1@radixweb.com"
$ws.Range("C2").Value = "Path"
$ws.Range("D2").Value = "9Prints_Web"
$ws.Range("E2").Value = "9Print_Mains"
$ws.Range("F2").Value = "9Prints"
$ws.Range("G2").Value = "Radixweb@8"
$ws.Range("H2").Value = 8789878978
$ws.Range("I2").Value = "Order Wise Sales Agent [ order.sales ] "
$ws.Range("J2").Value = 1
$ws.Range("K2").Value = "Discount"
$ws.Range("L2").Value = "Custom"
$ws.Range("M2").Value = 10
$ws.Range("N2").Value = "Store wise"
$ws.Range("O2").Value = 100000
$ws.Range("P2").Value = 0
$ws.Range("Q2").Value = 1
$ws.Range("R2").Value = 0
$ws.Range("S2").Value = 1
$ws.Range("T2").Value = "No"
$ws.Range("U2").Value = 0
$ws.Range("V2").Value = 0
$ws.Range("W2").Value = "No"
$ws.Range("X2").Value = 1
$ws.Range("Y2").Value = "Both"
$ws.Range("Z2").Value = 0
$ws.Range("AA2").Value = 1
$ws.Range("AB2").Value = 0
$ws.Range("A3").Value = "10Prints"
$ws.Range("B3").Value = "vaibhav.darji+store2@radixweb.com"
$ws.Range("C3").Value = "Subdomain"
$ws.Range("D3").Value = "10Prints_Web"
$ws.Range("E3").Value = "10Print_Mains"
$ws.Range("F3").Value = "10Prints"
$ws.Range("G3").Value = "Radixweb@9"
$ws.Range("H3").Value = 8789878979
$ws.Range("I3").Value = "General Sales Agent [ general.sales ] "
$ws.Range("J3").Value = 0
$ws.Range("K3").Value = "Markup"
$ws.Range("L3").Value = "Finesse"
$ws.Range("M3").Value = 20
$ws.Range("N3").Value = "Disable"
$ws.Range("O3").Value = 250000
$ws.Range("P3").Value = 1
$ws.Range("Q3").Value = 0
$ws.Range("R3").Value = 1
$ws.Range("S3").Value = 0
$ws.Range("T3").Value = "Yes"
$ws.Range("U3").Value = 1
$ws.Range("V3").Value = 1
$ws.Range("W3").Value = "Optional"
$ws.Range("X3").Value = 0
$ws.Range("Y3").Value = "Order wise only"
$ws.Range("Z3").Value = 1
$ws.Range("AA3").Value = 0
$ws.Range("AB3").Value = 1
$ws.Range("A4").Value = "11Prints"
$ws.Range("B4").Value = "vaibhav.darji+store3@radixweb.com"
$ws.Range("C4").Value = "Path"
$ws.Range("D4").Value = "11Prints_Web"
$ws.Range("E4").Value = "11Print_Mains"
$ws.Range("F4").Value = "11Prints"
$ws.Range("G4").Value = "Radixweb@10"
$ws.Range("H4").Value = 8789878980
$ws.Range("I4").Value = "kt sales [ kt sales ] "
$ws.Range("J4").Value = 1
$ws.Range("K4").Value = "Discount"
$ws.Range("L4").Value = "Crystals"
$ws.Range("M4").Value = 30
$ws.Range("N4").Value = "Customer wise"
$ws.Range("O4").Value = 500000
$ws.Range("P4").Value = 0
$ws.Range("Q4").Value = 1
$ws.Range("R4").Value = 0
$ws.Range("S4").Value = 1
$ws.Range("T4").Value = "Partial Order"
$ws.Range("U4").Value = 0
$ws.Range("V4").Value = 0
$ws.Range("W4").Value = "Mandatory"
$ws.Range("X4").Value = 1
$ws.Range("Y4").Value = "Stores"
$ws.Range("Z4").Value = 0
$ws.Range("AA4").Value = 1
$ws.Range("AB4").Value = 0

[void]$ws.Range("P6").Select()
